$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> D (Price) and/or E (Volume(1h)) new values,
# as produced by the GitHub Actions symbol-list refresh
# (commit: "Updated symbol list on Mon Jan 30 11:48:08 UTC 2023 with GitHub Actions").
# All Price/Volume columns in this sheet are stored as plain text, so values
# are written with a leading apostrophe to force text and the style is reset
# back to Normal afterwards to avoid leaving stray numeric formatting behind.
$updates = @(
    @{ Row = 2;  D = "308.38";     E = "-2.29%" },
    @{ Row = 3;  D = "38.04";      E = "-3.67%" },
    @{ Row = 4;  D = "5.055";      E = "-1.68%" },
    @{ Row = 5;  D = "0.07898";    E = "-3.42%" },
    @{ Row = 6;  D = "1.997";      E = "1.68%" },
    @{ Row = 7;  D = "4.370";      E = "3.13%" },
    @{ Row = 8;  D = "8.210";      E = "-0.21%" },
    @{ Row = 9;  D = $null;        E = "1.26%" },
    @{ Row = 10; D = "0.9254";     E = "-0.30%" },
    @{ Row = 11; D = "0.1280";     E = "-9.25%" },
    @{ Row = 12; D = "0.1877";     E = "-5.85%" },
    @{ Row = 13; D = "0.08678";    E = "-3.74%" },
    @{ Row = 14; D = "0.03462";    E = "-1.09%" },
    @{ Row = 15; D = "0.09729";    E = "-0.99%" },
    @{ Row = 16; D = $null;        E = "-0.62%" },
    @{ Row = 17; D = "0.006113";   E = "4.32%" },
    @{ Row = 18; D = "3.576";      E = "-2.11%" },
    @{ Row = 19; D = "0.3439";     E = "-0.72%" },
    @{ Row = 20; D = "0.1286";     E = "-1.29%" },
    @{ Row = 21; D = $null;        E = "5.20%" },
    @{ Row = 22; D = "0.2519";     E = "3.75%" },
    @{ Row = 23; D = "0.04338";    E = "-0.94%" },
    @{ Row = 24; D = "0.001220";   E = "-0.20%" },
    @{ Row = 25; D = $null;        E = "-4.01%" },
    @{ Row = 26; D = $null;        E = "176.35%" },
    @{ Row = 39; D = "0.02256";    E = "2.48%" },
    @{ Row = 40; D = "0.05018";    E = "-2.99%" },
    @{ Row = 41; D = "0.007530";   E = "-0.81%" },
    @{ Row = 42; D = "0.009891";   E = "1.50%" },
    @{ Row = 43; D = "0.1361";     E = "-1.07%" },
    @{ Row = 44; D = "0.002024";   E = "-4.99%" },
    @{ Row = 45; D = "0.008536";   E = "-6.49%" },
    @{ Row = 46; D = "0.00006435"; E = "0.50%" },
    @{ Row = 47; D = $null;        E = "0.19%" },
    @{ Row = 48; D = "0.003003";   E = "8.58%" },
    @{ Row = 49; D = "0.001204";   E = "0.34%" },
    @{ Row = 50; D = "0.00002104"; E = "0.19%" },
    @{ Row = 51; D = "0.0002004";  E = "0.19%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.Style = "Normal"
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.Style = "Normal"
        $cell.Value = "'" + $u.E
        $cell.Style = "Normal"
    }
}
